$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# --- New columns: X (Current VL balance), Y (Current SL balance), Z (ATM No./Account No.) ---

# 1. Give the brand-new, now-last column (Z) the "last column" look that W used to have
#    (bold header font, no fill) before W stops being the last column.
$ws.Range("W1").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)

# 2. Fill in the new header captions (shared strings get appended in order).
$ws.Range("X1").Value = "Current VL balance (hours)"
$ws.Range("Y1").Value = "Current SL balance (hours)"
$ws.Range("Z1").Value = "ATM No./Account No."

# 3. W is no longer the last column, so it now gets the regular highlighted header
#    look shared by the other inner header cells (same as V).
$ws.Range("V1").Copy()
$ws.Range("W1:Y1").PasteSpecial($xlPasteFormats)

# 4. Give the new columns sensible widths matching the new header text.
$ws.Columns.Item(24).ColumnWidth = 25.28
$ws.Columns.Item(25).ColumnWidth = 26.12
$ws.Columns.Item(26).ColumnWidth = 21.54

# --- Misc form clean-up ---

# 5. "Job position" header (L1) loses its highlight fill, matching the plain
#    header look used by most other columns (same as D1).
$ws.Range("D1").Copy()
$ws.Range("L1").PasteSpecial($xlPasteFormats)

# 6. "SSS No." header (N1) becomes center-aligned.
$ws.Range("N1").HorizontalAlignment = $xlCenter

$excel.CutCopyMode = $false

# 7. Update the remembered selection on the sheet.
$null = $ws.Range("C25").Select()
